# Update "想去人数" (want-to-go count) figures in column F across sheets,
# matching the latest data pull (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 7894
$ws1.Range("F6").Value  = 4908
$ws1.Range("F9").Value  = 1537
$ws1.Range("F13").Value = 1190
$ws1.Range("F20").Value = 1249
$ws1.Range("F24").Value = 1273
$ws1.Range("F33").Value = 17
$ws1.Range("F40").Value = 92
$ws1.Range("F43").Value = 440

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value  = 866
$ws3.Range("F9").Value  = 1778
$ws3.Range("F10").Value = 2675

# 全部类型 (All types - aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 866
$ws4.Range("F7").Value  = 7894
$ws4.Range("F9").Value  = 4908
$ws4.Range("F12").Value = 1537
$ws4.Range("F16").Value = 1778
$ws4.Range("F17").Value = 2675
$ws4.Range("F19").Value = 1190
$ws4.Range("F24").Value = 1249
$ws4.Range("F27").Value = 1273
$ws4.Range("F41").Value = 92
$ws4.Range("F44").Value = 440
